$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their original text formatting
# (values such as "0.9997" or "1.001" must remain text, not be converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.623.64'
$ws.Range('E2').Value = '  -4.52%  '
$ws.Range('D3').Value = '1.844.79'
$ws.Range('E3').Value = '  -3.93%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').Value = '313.04'
$ws.Range('E5').Value = '  -3.53%  '
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '0.4249'
$ws.Range('E7').Value = '  -6.96%  '
$ws.Range('E8').Value = '  -4.47%  '
$ws.Range('D9').Value = '43.77'
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').Value = '0.07215'
$ws.Range('E10').Value = '  -6.79%  '
$ws.Range('D11').Value = '0.8977'
$ws.Range('E11').Value = '  -7.98%  '
$ws.Range('E12').Value = '  -7.23%  '
$ws.Range('D13').Value = '1.829.80'
$ws.Range('E13').Value = '  -4.17%  '
$ws.Range('D14').Value = '6.572'
$ws.Range('E14').Value = '  -5.45%  '
$ws.Range('D15').Value = '5.334'
$ws.Range('E15').Value = '  -6.34%  '
$ws.Range('D16').Value = '0.06801'
$ws.Range('E16').Value = '  -2.63%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '77.48'
$ws.Range('E18').Value = '  -8.24%  '
$ws.Range('D19').Value = '0.000008851'
$ws.Range('E19').Value = '  -6.49%  '
$ws.Range('D20').Value = '0.9992'
$ws.Range('D21').Value = '15.36'
$ws.Range('D22').Value = '27.596.74'
$ws.Range('E22').Value = '  -4.66%  '
$ws.Range('D23').Value = '4.944'
$ws.Range('E24').Value = '  -2.85%  '
$ws.Range('D25').Value = '2.055.49'
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('D26').Value = '2.047'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').Value = '152.18'
$ws.Range('E27').Value = '  -3.55%  '
$ws.Range('D28').Value = '18.19'
$ws.Range('E28').Value = '  -4.45%  '
$ws.Range('D29').Value = '5.314'
$ws.Range('E29').Value = '  -5.14%  '
$ws.Range('D30').Value = '111.18'
$ws.Range('E30').Value = '  -5.59%  '
$ws.Range('D31').Value = '1.749'
$ws.Range('E31').Value = '  -4.81%  '
$ws.Range('D32').Value = '0.08885'
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('D33').Value = '0.7769'
$ws.Range('E33').Value = '  -10.01%  '
$ws.Range('D34').Value = '4.487'
$ws.Range('E34').Value = '  -11.99%  '
$ws.Range('D35').Value = '2.840'
$ws.Range('E35').Value = '  -5.66%  '
$ws.Range('D36').Value = '1.086'
$ws.Range('E36').Value = '  -12.32%  '
$ws.Range('D37').Value = '0.9992'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('D38').Value = '0.05428'
$ws.Range('E38').Value = '  -4.50%  '
$ws.Range('D39').Value = '1.085'
$ws.Range('E39').Value = '  -5.66%  '
$ws.Range('D40').Value = '2.979'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').Value = '0.01924'
$ws.Range('E41').Value = '  -5.88%  '
$ws.Range('E42').Value = '  -8.04%  '
$ws.Range('D43').Value = '6.777'
$ws.Range('E43').Value = '  -8.99%  '
$ws.Range('E44').Value = '  -6.82%  '
$ws.Range('D45').Value = '0.06616'
$ws.Range('E45').Value = '  -4.67%  '
$ws.Range('D46').Value = '8.228'
$ws.Range('E46').Value = '  -11.61%  '
$ws.Range('D47').Value = '106.22'
$ws.Range('E47').Value = '  -3.89%  '
$ws.Range('D48').Value = '0.4710'
$ws.Range('E48').Value = '  -8.59%  '
$ws.Range('D49').Value = '10.22'
$ws.Range('E49').Value = '  -8.99%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('E51').Value = '  -13.13%  '
